$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# D7: "高孟芳菲" -> rich text "（单）高孟芳菲\n（双）张笑微"
# with the "（单）" prefix carrying a distinct (CJK) font variant, and the
# remainder explicitly carrying the regular font. Build the two-run rich
# text by toggling the cell's base font around each Characters() call so the
# run captures the desired family.
$cell = $ws.Range("D7")
$cell.Value = "（单）高孟芳菲" + [char]10 + "（双）张笑微"
$cell.WrapText = $true

$cell.Font.Family = 2
$run2 = $cell.Characters(4, 11)
$run2.Font.Name = "等线"

$cell.Font.Family = 3
$run1 = $cell.Characters(1, 3)
$run1.Font.Name = "等线"

$ws.Rows(7).RowHeight = 27.75

# G9/H9: shift names - "李天元" drops off the roster, "时艺宁" moves into
# G9, and new hire "薛奕" takes H9.
$ws.Range("G9").Value = "时艺宁"
$ws.Range("H9").Value = "薛奕"

# Selection moves to E10 as last recorded in the file.
$ws.Range("E10").Select()
